$d = $word.ActiveDocument

# Find the (unique) old text and overwrite the Range's .Text directly —
# unlike Find.Execute(..., Replace:=wdReplaceAll) this does not run the
# 'smart quotes' autocorrect pass, so straight apostrophes are preserved.
function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { throw "Text not found: $old" }
    $rng.Text = $new
}

function Insert-ParasAfter($anchorText, [string[]]$newParas) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { throw "Anchor not found: $anchorText" }
    foreach ($p in $newParas) {
        $rng.InsertParagraphAfter()
        $rng = $d.Range($rng.End + 1, $rng.End + 1)
        $rng.InsertAfter($p)
    }
}

# --- Simple bullet text replacements ---
Replace-Text "• Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations" "• Lead comprehensive polling and research studies for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions"
Replace-Text "• Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics" "• Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing millions of records with millions of columns for electoral analytics and demographic analysis"
Replace-Text "• Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets" "• Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets"
Replace-Text "• Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering" "• Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering"
Replace-Text "• Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications" "• Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications"
Replace-Text "• Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices" "• Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices"
Replace-Text "• Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES" "• Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES"
Replace-Text "• Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions" "• Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions"
Replace-Text "• Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI" "• Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI"
Replace-Text "• Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products" "• Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company's distinguishing products"
Replace-Text "• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices" "• Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices"
Replace-Text "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research" "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions"

# --- Block replacements (rewrite first 4, then insert any extra) ---
Replace-Text "• Managed critical research operations for political campaigns" "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls"
Replace-Text "• Conducted comprehensive polling and demographic analysis" "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"
Replace-Text "• Developed strategic recommendations based on data analysis" "• Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"
Replace-Text "• Led research team in support of progressive political initiatives" "• Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs"
Insert-ParasAfter "• Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs" @("• Managed comprehensive research operations for progressive political initiatives and candidates")

Replace-Text "• Developed data analysis tools for political polling and research" "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"
Replace-Text "• Built statistical models for voter behavior analysis" "• Developed system that later became the Polling Consortium Database at The Analyst Institute"
Replace-Text "• Created data visualization tools for research presentations" "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions"
Replace-Text "• Supported senior researchers with technical analysis and reporting" "• Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle"
Insert-ParasAfter "• Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle" @("• Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps", "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding")

Replace-Text "• Managed field operations for political campaigns and research projects" "• Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions"
Replace-Text "• Developed data collection and management systems for field work" "• Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm"
Replace-Text "• Trained field staff on data collection protocols and quality control" "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"
Replace-Text "• Analyzed field data to inform campaign strategy and research findings" "• Created custom reports and data visualizations based on specific client requirements"
Insert-ParasAfter "• Created custom reports and data visualizations based on specific client requirements" @("• Introduced mapping and geospatial analysis into standard reporting procedures", "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL")

# --- Section heading renames (applies to both occurrences) ---
Replace-Text "Political Research and Data Analysis" "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns"
Replace-Text "Political Research and Data Analysis" "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns"
Replace-Text "Political Field Operations and Data Management" "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns"

